$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the typo in the email address shown in C4 (shared string reused by the mailto hyperlink)
$ws.Range("C4").Value = "acd@gmail.com"

# Fill in the missing "Sign Up" value in C5 (same text already used in D5)
$ws.Range("C5").Value = "Sign Up"

# Move the selection to C5, matching where the user left off editing
$ws.Range("C5").Select()
